$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (dated 2022-08-12 / serial 44785) was added to the
# "Rabanito" sheet. It slots in right after the existing row 244 (same
# market/category/quality/price/origin data, only the date differs), which
# pushes every subsequent record down by one row (old 245-304 -> new 246-305).

# Insert a new row at position 245, shifting rows 245:304 down to 246:305.
$ws.Rows.Item(245).Insert()

# Populate the new row 245 with the same data as row 244 (now still at 244),
# then overwrite just the date (column D) with the new record's date.
$ws.Range("A245:R245").Value = $ws.Range("A244:R244").Value2
$ws.Range("D245").Value = 44785
